{"js": "const pairs = [\n  [\"51\u00d725=\", \"20\u00d757=\"],\n  [\"67\u00d769=\", \"84\u00d781=\"],\n  [\"76\u00d729=\", \"79\u00d750=\"],\n  [\"28\u00d711=\", \"39\u00d776=\"],\n  [\"78\u00d723=\", \"68\u00d725=\"],\n  [\"71\u00d726=\", \"18\u00d722=\"],\n  [\"85\u00d787=\", \"21\u00d768=\"],\n  [\"27\u00d771=\", \"70\u00d761=\"],\n  [\"90\u00d790=\", \"89\u00d775=\"],\n  [\"65\u00d738=\", \"87\u00d732=\"],\n  [\"68\u00d763=\", \"68\u00d799=\"],\n  [\"93\u00d716=\", \"40\u00d780=\"],\n  [\"89\u00d787=\", \"60\u00d774=\"],\n  [\"90\u00d756=\", \"40\u00d739=\"],\n  [\"31\u00d727=\", \"40\u00d772=\"],\n  [\"57\u00d735=\", \"54\u00d733=\"],\n  [\"74\u00d794=\", \"95\u00d763=\"],\n  [\"15\u00d733=\", \"83\u00d768=\"],\n  [\"33\u00d773=\", \"44\u00d758=\"],\n  [\"69\u00d737=\", \"23\u00d763=\"],\n  [\"84\u00d727=\", \"62\u00d743=\"],\n  [\"24\u00d781=\", \"47\u00d717=\"],\n  [\"96\u00d780=\", \"79\u00d769=\"],\n  [\"85\u00d786=\", \"72\u00d788=\"],\n  [\"99\u00d755=\", \"98\u00d761=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [before, after] of pairs) {\n  const results = body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + before);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Before = \"51\u00d725=\"; After = \"20\u00d757=\" }\n    @{ Before = \"67\u00d769=\"; After = \"84\u00d781=\" }\n    @{ Before = \"76\u00d729=\"; After = \"79\u00d750=\" }\n    @{ Before = \"28\u00d711=\"; After = \"39\u00d776=\" }\n    @{ Before = \"78\u00d723=\"; After = \"68\u00d725=\" }\n    @{ Before = \"71\u00d726=\"; After = \"18\u00d722=\" }\n    @{ Before = \"85\u00d787=\"; After = \"21\u00d768=\" }\n    @{ Before = \"27\u00d771=\"; After = \"70\u00d761=\" }\n    @{ Before = \"90\u00d790=\"; After = \"89\u00d775=\" }\n    @{ Before = \"65\u00d738=\"; After = \"87\u00d732=\" }\n    @{ Before = \"68\u00d763=\"; After = \"68\u00d799=\" }\n    @{ Before = \"93\u00d716=\"; After = \"40\u00d780=\" }\n    @{ Before = \"89\u00d787=\"; After = \"60\u00d774=\" }\n    @{ Before = \"90\u00d756=\"; After = \"40\u00d739=\" }\n    @{ Before = \"31\u00d727=\"; After = \"40\u00d772=\" }\n    @{ Before = \"57\u00d735=\"; After = \"54\u00d733=\" }\n    @{ Before = \"74\u00d794=\"; After = \"95\u00d763=\" }\n    @{ Before = \"15\u00d733=\"; After = \"83\u00d768=\" }\n    @{ Before = \"33\u00d773=\"; After = \"44\u00d758=\" }\n    @{ Before = \"69\u00d737=\"; After = \"23\u00d763=\" }\n    @{ Before = \"84\u00d727=\"; After = \"62\u00d743=\" }\n    @{ Before = \"24\u00d781=\"; After = \"47\u00d717=\" }\n    @{ Before = \"96\u00d780=\"; After = \"79\u00d769=\" }\n    @{ Before = \"85\u00d786=\"; After = \"72\u00d788=\" }\n    @{ Before = \"99\u00d755=\"; After = \"98\u00d761=\" }\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Before\n    $find.Replacement.Text = $pair.After\n    $find.Forward = $true\n    $find.Wrap = $wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $found = $find.Execute($pair.Before, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $pair.After, $wdReplaceAll)\n\n    if (-not $found) {\n        throw \"No match found for: $($pair.Before)\"\n    }\n}\n"}
